# Escaleta MA_07_05_CO.xlsx — apply "Escaletas ubicadas en git hub" edit
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja2")

function Set-Revision($row, $dateSerial) {
    $ws.Range("V$row").Value = "Clara Melo"
    if ($dateSerial -ne $null) {
        $ws.Range("W$row").Value = $dateSerial
        $ws.Range("W$row").NumberFormat = "mmm-yy"
    }
}

# --- New column V ("responsable") + W (fecha revisión) data for the reviewed rows ---
Set-Revision 3 45261
Set-Revision 5 47088
Set-Revision 7 37987
Set-Revision 9 37987
Set-Revision 10 47088
Set-Revision 12 38353
Set-Revision 13 38353
Set-Revision 14 38353
Set-Revision 17 38718
Set-Revision 18 38718
Set-Revision 19 38718
Set-Revision 20 38718
Set-Revision 23 39814
Set-Revision 24 $null
Set-Revision 26 40909
Set-Revision 27 41275
Set-Revision 28 40909
Set-Revision 29 40909
Set-Revision 30 41640
Set-Revision 31 41640
Set-Revision 34 42005
Set-Revision 35 41640
Set-Revision 39 42005
Set-Revision 40 42005

# --- Hide columns that were hidden in this revision ---
$ws.Range("J1:L1").EntireColumn.Hidden = $true
$ws.Range("O1:P1").EntireColumn.Hidden = $true

# --- View state: scroll/selection/active sheet ---
$ws.Activate() | Out-Null
$ws.Range("X42").Select() | Out-Null

Write-Output "done"
